# Update column F ("dSF") values for several rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -4
$ws.Range("F7").Value  = -3
$ws.Range("F8").Value  = -1
$ws.Range("F9").Value  = -3
$ws.Range("F10").Value = -8
$ws.Range("F14").Value = 4
$ws.Range("F17").Value = 3
$ws.Range("F19").Value = -2
